# Update cryptos list figures (prices, 1h volume %, and a name/link swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.761.85'
$ws.Range('E2').Value = '  +4.44%  '
$ws.Range('D3').Value = '2.254.67'
$ws.Range('E3').Value = '  +4.01%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '70.65'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.97%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.661'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +16.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +9.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.58'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0964'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.56%  '
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '2.580.72'
$ws.Range('E15').Value = '  +3.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.880'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').Value = '2.249.36'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').Value = '42.673.39'
$ws.Range('E19').Value = '  +4.40%  '
$ws.Range('D20').Value = '0.0₃0989'
$ws.Range('E20').Value = '  +5.57%  '
$ws.Range('E21').Value = '  +3.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('E25').Value = '  +6.49%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.79%  '
$ws.Range('E34').Value = '  +5.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0794'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +25.05%  '
$ws.Range('E37').Value = '  +3.59%  '
$ws.Range('E38').Value = '  +10.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.71'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.54%  '
$ws.Range('E40').Value = '  +7.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.30'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.82%  '
$ws.Range('E42').Value = '  +6.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.33'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.05'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.51%  '
$ws.Range('E46').Value = '  +4.95%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('E48').Value = '  +2.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('E51').Value = '  +4.24%  '
